# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.559.57"
$ws.Range("E2").Value = "  +10.46%  "

$ws.Range("D3").Value = "3.259.49"
$ws.Range("E3").Value = "  +6.38%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.88%  "

$ws.Range("E7").Value = "  +4.59%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.50%  "

$ws.Range("E12").Value = "  +2.33%  "

$ws.Range("D13").Value = "3.767.74"
$ws.Range("E13").Value = "  +6.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.17%  "

$ws.Range("D16").Value = "3.251.97"
$ws.Range("E16").Value = "  +6.27%  "

$ws.Range("E17").Value = "  +5.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.98%  "

$ws.Range("D19").Value = "56.458.10"
$ws.Range("E19").Value = "  +10.23%  "

$ws.Range("E20").Value = "  +4.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "299.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.33%  "

$ws.Range("E24").Value = "  +8.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.53%  "

$ws.Range("E28").Value = "  +4.19%  "

$ws.Range("E29").Value = "  +2.03%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  +6.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.66%  "

$ws.Range("E35").Value = "  +3.41%  "

$ws.Range("E36").Value = "  +2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +26.26%  "

$ws.Range("E39").Value = "  +5.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.36%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.60%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.94%  "

$ws.Range("E44").Value = "  +4.67%  "

$ws.Range("E45").Value = "  +7.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.287"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.89%  "

$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("E48").Value = "  +56.64%  "

$ws.Range("D49").Value = "2.151.83"
$ws.Range("E49").Value = "  +4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("E51").Value = "  -3.59%  "
